$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "id_kategori" header value from J1 (keep the cell/style, just clear content)
$ws.Range("J1").ClearContents() | Out-Null

# Remove the stray value in J2 (the id_kategori data cell) entirely
$ws.Range("J2").ClearContents() | Out-Null

# Update the active selection to J6 as recorded in the saved view state
$ws.Range("J6").Select() | Out-Null
